$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3428084429008038
$ws.Range("C2").Value = 0.4025838691005106
$ws.Range("D2").Value = 0.07676214029817174
$ws.Range("E2").Value = 0.1109755168374491
$ws.Range("G2").Value = 0.002572678534647335
$ws.Range("I2").Value = 4.687047288554481
$ws.Range("K2").Value = 0.7055059354176478
$ws.Range("M2").Value = 0.3135549174527412
$ws.Range("B3").Value = 0.3369770500565892
$ws.Range("C3").Value = 0.3730551800777278
$ws.Range("D3").Value = 0.06987585593830659
$ws.Range("E3").Value = 0.1004564779350261
$ws.Range("G3").Value = 0.002578559251130494
$ws.Range("I3").Value = 4.346998500202233
$ws.Range("K3").Value = 0.6753379242549045
$ws.Range("M3").Value = 0.2929445366650967
$ws.Range("B4").Value = 0.3340559129126035
$ws.Range("C4").Value = 0.355063877279548
$ws.Range("D4").Value = 0.06569142156747887
$ws.Range("E4").Value = 0.09404779429636534
$ws.Range("G4").Value = 0.002582348256191192
$ws.Range("I4").Value = 4.137664757909022
$ws.Range("K4").Value = 0.6576085255159967
$ws.Range("M4").Value = 0.280549526496813
$ws.Range("B5").Value = 0.3330297943573015
$ws.Range("C5").Value = 0.3477666149499044
$ws.Range("D5").Value = 0.06399703412063218
$ws.Range("E5").Value = 0.09144842194680081
$ws.Range("G5").Value = 0.002583937303702359
$ws.Range("I5").Value = 4.052212590666244
$ws.Range("K5").Value = 0.6505814135778962
$ws.Range("M5").Value = 0.2755628918220623
$ws.Range("B6").Value = 0.3328692851907533
$ws.Range("C6").Value = 0.3465569633296468
$ws.Range("D6").Value = 0.06371632901641533
$ws.Range("E6").Value = 0.09101752548487241
$ws.Range("G6").Value = 0.002584203886853864
$ws.Range("I6").Value = 4.038014140695253
$ws.Range("K6").Value = 0.6494264574575652
$ws.Range("M6").Value = 0.2747387302704141
$ws.Range("B7").Value = 0.3340414112062291
$ws.Range("C7").Value = 0.3549653257512659
$ws.Range("D7").Value = 0.06566852698918524
$ws.Range("E7").Value = 0.09401268926172435
$ws.Range("G7").Value = 0.002582369504327559
$ws.Range("I7").Value = 4.136512929055698
$ws.Range("K7").Value = 0.6575129571321554
$ws.Range("M7").Value = 0.280482015300187
$ws.Range("B8").Value = 0.3406599940489201
$ws.Range("C8").Value = 0.3923729907542679
$ws.Range("D8").Value = 0.07437853898568392
$ws.Range("E8").Value = 0.1073379391361442
$ws.Range("G8").Value = 0.002574669326024547
$ws.Range("I8").Value = 4.56990454776664
$ws.Range("K8").Value = 0.6949380896863886
$ws.Range("M8").Value = 0.3063939690205757
$ws.Range("B9").Value = 0.358942652844064
$ws.Range("C9").Value = 0.4668724960314989
$ws.Range("D9").Value = 0.09181636502879087
$ws.Range("E9").Value = 0.1338838641308087
$ws.Range("G9").Value = 0.00256097504593471
$ws.Range("I9").Value = 5.415999429335443
$ws.Range("K9").Value = 0.7747183933206259
$ws.Range("M9").Value = 0.3593134431116312
$ws.Range("B10").Value = 0.375711746789591
$ws.Range("C10").Value = 0.522361489121522
$ws.Range("D10").Value = 0.1048615259798567
$ws.Range("E10").Value = 0.1536673728393083
$ws.Range("G10").Value = 0.002551759061235396
$ws.Range("I10").Value = 6.03609766053853
$ws.Range("K10").Value = 0.8373632523684478
$ws.Range("M10").Value = 0.3995439904804314
$ws.Range("B11").Value = 0.3840881047399023
$ws.Range("C11").Value = 0.5477816231671682
$ws.Range("D11").Value = 0.1108503612981622
$ws.Range("E11").Value = 0.1627343312406211
$ws.Range("G11").Value = 0.00254774746416999
$ws.Range("I11").Value = 6.318039103252829
$ws.Range("K11").Value = 0.8667675288001533
$ws.Range("M11").Value = 0.4181541572932588
$ws.Range("B12").Value = 0.3873694648953006
$ws.Range("C12").Value = 0.5574341117093127
$ws.Range("D12").Value = 0.1131262842272633
$ws.Range("E12").Value = 0.1661779000498314
$ws.Range("G12").Value = 0.002546254179723752
$ws.Range("I12").Value = 6.424795493171075
$ws.Range("K12").Value = 0.878034955018876
$ws.Range("M12").Value = 0.4252469508129764
$ws.Range("B13").Value = 0.3866578690823417
$ws.Range("C13").Value = 0.5553540860911994
$ws.Range("D13").Value = 0.1126357602170316
$ws.Range("E13").Value = 0.1654358094628137
$ws.Range("G13").Value = 0.00254657463976867
$ws.Range("I13").Value = 6.401803787892334
$ws.Range("K13").Value = 0.8756023762107361
$ws.Range("M13").Value = 0.4237173485184655
$ws.Range("B14").Value = 0.3843558606209285
$ws.Range("C14").Value = 0.5485752044487526
$ws.Range("D14").Value = 0.111037439496485
$ws.Range("E14").Value = 0.1630174306660592
$ws.Range("G14").Value = 0.002547624094217803
$ws.Range("I14").Value = 6.326822120567044
$ws.Range("K14").Value = 0.8676918347605351
$ws.Range("M14").Value = 0.4187367669732467
$ws.Range("B15").Value = 0.3829601187035507
$ws.Range("C15").Value = 0.5444264120208118
$ws.Range("D15").Value = 0.1100594815107456
$ws.Range("E15").Value = 0.1615374329987063
$ws.Range("G15").Value = 0.002548270272883181
$ws.Range("I15").Value = 6.280892910624345
$ws.Range("K15").Value = 0.8628637447767176
$ws.Range("M15").Value = 0.4156919784486064
$ws.Range("B16").Value = 0.3751795414170829
$ws.Range("C16").Value = 0.5207038759892839
$ws.Range("D16").Value = 0.1044712591179007
$ws.Range("E16").Value = 0.1530762169850206
$ws.Range("G16").Value = 0.002552024851408854
$ws.Range("I16").Value = 6.017670243611093
$ws.Range("K16").Value = 0.8354600536592329
$ws.Range("M16").Value = 0.3983340837916387
$ws.Range("B17").Value = 0.3705992983081501
$ws.Range("C17").Value = 0.5061970358735834
$ws.Range("D17").Value = 0.1010572002795698
$ws.Range("E17").Value = 0.1479030823171215
$ws.Range("G17").Value = 0.002554374344281834
$ws.Range("I17").Value = 5.856163186903217
$ws.Range("K17").Value = 0.8188826291386135
$ws.Range("M17").Value = 0.3877655039485504
$ws.Range("B18").Value = 0.3680351569938409
$ws.Range("C18").Value = 0.4978697789334205
$ws.Range("D18").Value = 0.0990986414959707
$ws.Range("E18").Value = 0.1449339446177049
$ws.Range("G18").Value = 0.002555742738790556
$ws.Range("I18").Value = 5.763254278335495
$ws.Range("K18").Value = 0.8094329480325086
$ws.Range("M18").Value = 0.3817157631771977
$ws.Range("B19").Value = 0.3671789999763746
$ws.Range("C19").Value = 0.4950531543778425
$ws.Range("D19").Value = 0.09843637886405077
$ws.Range("E19").Value = 0.1439297169890779
$ws.Range("G19").Value = 0.002556208983901861
$ws.Range("I19").Value = 5.731794095457758
$ws.Range("K19").Value = 0.8062480159774736
$ws.Range("M19").Value = 0.3796723777594337
$ws.Range("B20").Value = 0.3710795846238284
$ws.Range("C20").Value = 0.5077395799179385
$ws.Range("D20").Value = 0.1014201013010165
$ws.Range("E20").Value = 0.1484531149593238
$ws.Range("G20").Value = 0.002554122475549994
$ws.Range("I20").Value = 5.873357284492442
$ws.Range("K20").Value = 0.8206384886090348
$ws.Range("M20").Value = 0.3888875353987444
$ws.Range("B21").Value = 0.3850290318513885
$ws.Range("C21").Value = 0.5505656019679463
$ws.Range("D21").Value = 0.1115066837796803
$ws.Range("E21").Value = 0.1637274893725618
$ws.Range("G21").Value = 0.002547315144847941
$ws.Range("I21").Value = 6.348846193260215
$ws.Range("K21").Value = 0.8700117337219808
$ws.Range("M21").Value = 0.4201984401742465
$ws.Range("B22").Value = 0.3947844458495524
$ws.Range("C22").Value = 0.5787093793788358
$ws.Range("D22").Value = 0.1181460660639004
$ws.Range("E22").Value = 0.1737692991173674
$ws.Range("G22").Value = 0.002543016578669899
$ws.Range("I22").Value = 6.659560636549486
$ws.Range("K22").Value = 0.903054366896697
$ws.Range("M22").Value = 0.4409278132319514
$ws.Range("B23").Value = 0.3895187394163315
$ws.Range("C23").Value = 0.563674092438589
$ws.Range("D23").Value = 0.1145981015604747
$ws.Range("E23").Value = 0.1684042438789604
$ws.Range("G23").Value = 0.002545297099904363
$ws.Range("I23").Value = 6.493726572825153
$ws.Range("K23").Value = 0.885347274439539
$ws.Range("M23").Value = 0.4298394645742292
$ws.Range("B24").Value = 0.3708622320845905
$ws.Range("C24").Value = 0.507042155581189
$ws.Range("D24").Value = 0.1012560205942634
$ws.Range("E24").Value = 0.1482044297648812
$ws.Range("G24").Value = 0.002554236290380429
$ws.Range("I24").Value = 5.865584009398077
$ws.Range("K24").Value = 0.8198444127242794
$ws.Range("M24").Value = 0.3883801832159008
$ws.Range("B25").Value = 0.3534183430525673
$ws.Range("C25").Value = 0.4465901665128627
$ws.Range("D25").Value = 0.08705909837364345
$ws.Range("E25").Value = 0.126655266782187
$ws.Range("G25").Value = 0.002564530434037593
$ws.Range("I25").Value = 5.187433704547914
$ws.Range("K25").Value = 0.7524378768213751
$ws.Range("M25").Value = 0.3447650345438902
